# Assembly in real world is done and some Joints CAD are modified.
# All BOM components in column G ("Status") are now confirmed "Arrived" -
# including the two rows that previously flagged problems
# ("wrong Part + Bolts???" / "NOT arrived") and the rows that had not
# been marked yet. Every status cell also gets the same bordered,
# centered "Arrived" look (matching the style already used on most of
# the column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2 already carries the bordered/centered "Arrived" style - use it as
# the formatting template for the rest of the Status column (G3:G21).
$ws.Range("G2").Copy()
for ($r = 3; $r -le 21; $r++) {
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)
}

# Now mark every component in the BOM as Arrived.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 7).Value = "Arrived"
}
